$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header for row 2 (B2 now "Fase 1" instead of "Tasca")
$ws.Range("B2").Value = "Fase 1"

# New section header "Fase 2" at B11
$ws.Range("B11").Value = "Fase 2"

# New rows 13-16
$ws.Range("B13").Value = "Opció Connect"
$ws.Range("C13").Value = 15

$ws.Range("B14").Value = "Opció Show Connections"
$ws.Range("C14").Value = 7

$ws.Range("B15").Value = "Opció Say"
$ws.Range("C15").Value = 1

$ws.Range("B16").Value = "Opció Exit"
$ws.Range("C16").Value = 40

$ws.Range("C13").Select()
